$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.23260486125946
$ws.Range("B1").Value = 2.519661664962769
$ws.Range("C1").Value = 4.509186744689941
$ws.Range("D1").Value = 2.492693662643433
$ws.Range("E1").Value = 1.069980502128601
